$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B16").Value = "I am still not seeing April in the solutions."
$ws.Range("B16").Select()
